$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8538402915000916
$ws.Range("B1").Value = 1.084203720092773
$ws.Range("C1").Value = 1.599539756774902
$ws.Range("D1").Value = 2.268487691879272
$ws.Range("E1").Value = 1.74068295955658
